# "added a new exclusion template"
#
# The KPI exclusions template's generic "Empty; Irrelevant" exclusion value
# is expanded to also cover "General Empty", the active selection is moved,
# and the data columns are nudged slightly wider (cosmetic re-layout that
# happened alongside the content edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Slightly widen the tab bar relative to the horizontal scrollbar
# (workbookView tabRatio 990 -> 993 in the saved view state).
$wb.Windows.Item(1).TabRatio = 0.993

# Update both cells that shared the "Empty; Irrelevant" exclusion text so the
# underlying shared string is updated in place rather than forking a new one.
$newValue = "Empty; Irrelevant; General Empty"
$ws.Range("C2").Value = $newValue
$ws.Range("C3").Value = $newValue

# Widen the first three data columns slightly.
$ws.Columns.Item(1).ColumnWidth = 25.7085020242915
$ws.Columns.Item(2).ColumnWidth = 18.8542510121457
$ws.Columns.Item(3).ColumnWidth = 44.8825910931174

# Move the active selection to C22, as captured in the saved view state.
$ws.Range("C22").Select()
